$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget Summary")
$ws.Activate()
Write-Host ("A1 value: " + $ws.Range("A1").Value)
